# Applies the cryptos.xlsx price/volume refresh described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new value map (applied in row order, matching the source diff).
$edits = [ordered]@{
    'D2' = '66.146.80'
    'E2' = '  -0.56%  '
    'D3' = '3.334.24'
    'E3' = '  -0.35%  '
    'D4' = '0.999'
    'E4' = '  -0.03%  '
    'D5' = '584.88'
    'E5' = '  +3.22%  '
    'D6' = '185.31'
    'E6' = '  -2.91%  '
    'E7' = '  -0.08%  '
    'D8' = '3.329.57'
    'E8' = '  +0.10%  '
    'E9' = '  -2.41%  '
    'E10' = '  -2.99%  '
    'E11' = '  -2.04%  '
    'D12' = '47.10'
    'E12' = '  -2.34%  '
    'E13' = '  -2.12%  '
    'D14' = '677.57'
    'E14' = '  +11.55%  '
    'D15' = '3.864.98'
    'E15' = '  -0.13%  '
    'D16' = '8.50'
    'E16' = '  -2.99%  '
    'D17' = '66.306.68'
    'E17' = '  -0.26%  '
    'B18' = 'TRON'
    'C18' = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
    'D18' = '0.118'
    'E18' = '  -0.69%  '
    'B19' = 'Chainlink'
    'C19' = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
    'D19' = '17.91'
    'E19' = '  -1.55%  '
    'D20' = '3.336.22'
    'E20' = '  -0.23%  '
    'E21' = '  -1.08%  '
    'E22' = '  -2.39%  '
    'D23' = '17.76'
    'E23' = '  -5.29%  '
    'D24' = '103.73'
    'E24' = '  +2.56%  '
    'D25' = '5.03'
    'E25' = '  -3.01%  '
    'E26' = '  -1.52%  '
    'E27' = '  +0.18%  '
    'D28' = '9.44'
    'E28' = '  -3.91%  '
    'D29' = '32.51'
    'E29' = '  +6.00%  '
    'E30' = '  -3.03%  '
    'D31' = '6.82'
    'E31' = '  -0.62%  '
    'D32' = '609.32'
    'E32' = '  +6.59%  '
    'D33' = '3.90'
    'E33' = '  -3.24%  '
    'D34' = '11.11'
    'E34' = '  -0.87%  '
    'D35' = '3.847.83'
    'E35' = '  +3.56%  '
    'E36' = '  -1.36%  '
    'E37' = '  +0.09%  '
    'D38' = '56.05'
    'E38' = '  -2.20%  '
    'D39' = '0.128'
    'E39' = '  -2.49%  '
    'E40' = '  -1.83%  '
    'D41' = '0.0₃0703'
    'E41' = '  -4.65%  '
    'D42' = '3.19'
    'E42' = '  -4.02%  '
    'D43' = '32.63'
    'E43' = '  -5.00%  '
    'D44' = '3.40'
    'E44' = '  +3.48%  '
    'E45' = '  -2.98%  '
    'E46' = '  -3.12%  '
    'E47' = '  -12.53%  '
    'E48' = '  -1.88%  '
    'B49' = 'FirstDigitalUSD'
    'C49' = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
    'D49' = '1.00'
    'E49' = '  +0.35%  '
    'B50' = 'ThetaToken'
    'C50' = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
    'D50' = '2.55'
    'E50' = '  -2.74%  '
    'E51' = '  +1.60%  '
}

# Cells whose new text is a bare number/decimal: Excel's COM layer would
# silently coerce these to numeric cells on assignment, but the source
# workbook stores every Price/Volume cell as text (inline string). Force
# text formatting first, then restore the default (unstyled) cell style so
# no stray number-format styling is left behind.
$numericLooking = @(
    'D4'
    'D5'
    'D6'
    'D12'
    'D14'
    'D16'
    'D18'
    'D19'
    'D23'
    'D24'
    'D25'
    'D28'
    'D29'
    'D31'
    'D32'
    'D33'
    'D34'
    'D38'
    'D39'
    'D42'
    'D43'
    'D44'
    'D49'
    'D50'
)
foreach ($cell in $numericLooking) {
    $ws.Range($cell).NumberFormat = "@"
}

foreach ($cell in $edits.Keys) {
    $ws.Range($cell).Value = $edits[$cell]
}

foreach ($cell in $numericLooking) {
    $ws.Range($cell).Style = "Normal"
}
